$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final state for data rows 8-17 (columns A:E).
# A = bus index, B = name, C/D = values, E = in_service flag.
$rows = @(
    @(8,  6,  "line7", 14, 11, $true),
    @(9,  7,  "line8", 16, 9,  $true),
    @(10, 8,  "extr1", 5,  12, $true),
    @(11, 9,  "extr2", 5,  9,  $true),
    @(12, 10, "extr3", 10, 11, $true),
    @(13, 11, "extr4", 7,  8,  $true),
    @(14, 12, "extr5", 9,  11, $false),
    @(15, 13, "extr6", 7,  11, $false),
    @(16, 14, "extr7", 5,  7,  $true),
    @(17, 15, "extr8", 8,  5,  $true)
)

foreach ($r in $rows) {
    $rowIdx = $r[0]
    $ws.Cells.Item($rowIdx, 1).Value = $r[1]
    $ws.Cells.Item($rowIdx, 2).Value = $r[2]
    $ws.Cells.Item($rowIdx, 3).Value = $r[3]
    $ws.Cells.Item($rowIdx, 4).Value = $r[4]
    $ws.Cells.Item($rowIdx, 5).Value = $r[5]
}

# Rows 16 & 17 are brand new - give column A the same formatting (bold,
# bordered, centered) already used by the rest of the A column.
$ws.Range("A15").Copy() | Out-Null
$ws.Range("A16:A17").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
